# Update the "Corr/total marks" figures in the marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row: points per right answer (B11) 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total score (B12) 57 -> 95
$ws.Range("B12").Value = 95

# "Total" row: correct/total marks label (E12) "55/84" -> "95/140"
$ws.Range("E12").Value = "95/140"
